$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 48, shifting existing rows 48:86 down to 51:89
$ws.Rows("48:50").Insert()

# Populate the 3 new rows (48, 49, 50) with new weekly price data
$newRows = @(
    @{ Row = 48; D = 44679; L = "Especial"; M = 240; N = 12000; O = 13000; P = 12500; S = 694 },
    @{ Row = 49; D = 44679; L = "Primera";  M = 300; N = 10000; O = 11000; P = 10500; S = 583 },
    @{ Row = 50; D = 44679; L = "Segunda";  M = 240; N = 8000;  O = 9000;  P = 8500;  S = 472 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100107
    $ws.Cells.Item($row, 8).Value = "Otros"
    $ws.Cells.Item($row, 9).Value = 100107011
    $ws.Cells.Item($row, 10).Value = "Tuna"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = "`$/caja 18 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 18
}
